$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3992870.79
$ws.Range("C7").Value = -10.13277997461813
$ws.Range("D7").Value = 3466
$ws.Range("E7").Value = 3466
$ws.Range("F7").Value = 1152.011191575303
$ws.Range("G7").Value = 22.79606290831175
